$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.958.57"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.406.04"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.91"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.58"
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  -1.70%  "

$ws.Range("E11").Value = "  -2.36%  "

$ws.Range("E12").Value = "  -0.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.61"
$ws.Range("E13").Value = "  -2.72%  "

$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.840.87"
$ws.Range("E15").Value = "  -0.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.954.64"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.403.90"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.23"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "322.20"
$ws.Range("E19").Value = "  -0.95%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.14"
$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.19"
$ws.Range("E23").Value = "  +1.89%  "

$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.80"
$ws.Range("E25").Value = "  -4.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "568.16"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.519.95"
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0937"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.19"
$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("E31").Value = "  -2.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.147"
$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("E34").Value = "  -3.15%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.69"
$ws.Range("E36").Value = "  -2.77%  "

$ws.Range("E37").Value = "  -5.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.379"
$ws.Range("E38").Value = "  -1.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "151.16"
$ws.Range("E39").Value = "  +2.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.61"
$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("E41").Value = "  -9.14%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.76"
$ws.Range("E44").Value = "  -3.07%  "

$ws.Range("E45").Value = "  -0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0531"
$ws.Range("E46").Value = "  -2.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.88"
$ws.Range("E47").Value = "  -2.89%  "

$ws.Range("E48").Value = "  -0.68%  "

$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.54"
$ws.Range("E51").Value = "  +0.48%  "
